# songs.xlsx — "added popularity of track. Added functionality to split
# attributes into excel columns"
#
# The track that used to duplicate "Too Good - Drake" in row 6 is replaced
# with a new track, and column A is widened so the longer titles still fit.
# A block of narrower columns to the right of A is reserved for the new
# per-track attribute columns (e.g. popularity) that the split introduces.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the duplicated "Too Good - Drake" entry in A6 with the new track.
$ws.Range("A6").Value = "Trust Issues Drake"

# Column A needs to be a bit wider to comfortably fit the song titles.
$ws.Columns.Item(1).ColumnWidth = 32.1377551020408

# Reserve narrower columns to the right (B:AMK, i.e. up to column 1025) for
# the attributes (popularity, etc.) that get split out of column A.
$ws.Range("B1:AMK1").EntireColumn.ColumnWidth = 8.50510204081633

# Leave the active selection on the row that was just edited.
$ws.Range("A7").Select()
